$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels for the two new transformer parameter columns
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Update existing parameter values in row 2
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 90
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# New transformer parameter values
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1

# Update the active cell selection to match the saved workbook state
$ws.Range("F2").Select()
